$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91828
$ws.Range("B3").Value = 58520
$ws.Range("B4").Value = 57884
$ws.Range("B5").Value = 57884
